$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.143.84'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '3.056.91'
$ws.Range('E3').Value = '  -3.02%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''589.06'
$ws.Range('D6').Value = '''152.01'
$ws.Range('E6').Value = '  +4.00%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '''0.546'
$ws.Range('E8').Value = '  +3.02%  '
$ws.Range('D9').Value = '3.060.87'
$ws.Range('E9').Value = '  -2.61%  '
$ws.Range('E10').Value = '  -4.06%  '
$ws.Range('D11').Value = '''5.81'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').Value = '''0.0000239'
$ws.Range('E13').Value = '  -3.36%  '
$ws.Range('D14').Value = '''36.95'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').Value = '3.564.81'
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').Value = '''7.18'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').Value = '63.187.86'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '3.059.07'
$ws.Range('E19').Value = '  -2.82%  '
$ws.Range('D20').Value = '''473.61'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').Value = '''14.63'
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('D22').Value = '''0.715'
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').Value = '''7.52'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').Value = '''2.37'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('D25').Value = '''12.95'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('D26').Value = '''81.09'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('D28').Value = '''9.97'
$ws.Range('E28').Value = '  +2.09%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '''7.27'
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('D33').Value = '''0.113'
$ws.Range('E33').Value = '  +2.24%  '
$ws.Range('D34').Value = '''27.17'
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('D35').Value = '0.0₃0840'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E36').Value = '  -2.41%  '
$ws.Range('D37').Value = '''6.09'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('E39').Value = '  -4.97%  '
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = '''50.36'
$ws.Range('E41').Value = '  -2.08%  '
$ws.Range('D42').Value = '''443.36'
$ws.Range('E42').Value = '  -4.41%  '
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').Value = '''40.29'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0362'
$ws.Range('E45').Value = '  -2.82%  '
$ws.Range('E46').Value = '  +1.95%  '
$ws.Range('D47').Value = '2.793.05'
$ws.Range('E47').Value = '  -4.64%  '
$ws.Range('D48').Value = '''131.19'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '''25.03'
$ws.Range('E50').Value = '  +3.18%  '
$ws.Range('E51').Value = '  +0.18%  '

Write-Output "Applied 90 cell updates"
